# "Generate Report for handback" — mark the two handed-off files as handed
# back (in sync with en-US), populate the "Latest Target File" / "Latest
# Handback File" columns on the per-locale sheets, and stamp the handback
# datetime.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status column for both locales -----------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet ------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsZh.Hyperlinks.Add(
    $wsZh.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d0f4f5c5fc27c12fe0d02e1e3d4264163c74298/e2e/6a104238-6086-4051-b5c2-9444b42ccc29.md",
    $null,
    $null,
    "6a104238-6086-4051-b5c2-9444b42ccc29.md")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff2705860e7492cff4c6fa7ace9133cb47bc0bb9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6a104238-6086-4051-b5c2-9444b42ccc29.4a63ded2fe0ab905558bc9b658fa27f0c8174700.zh-cn.xlf",
    $null,
    $null,
    "6a104238-6086-4051-b5c2-9444b42ccc29.4a63ded2fe0ab905558bc9b658fa27f0c8174700.zh-cn.xlf")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d0f4f5c5fc27c12fe0d02e1e3d4264163c74298/e2e/878e5e59-2378-49c6-bcb0-0a60edbf870d.md",
    $null,
    $null,
    "878e5e59-2378-49c6-bcb0-0a60edbf870d.md")

$wsZh.Hyperlinks.Add(
    $wsZh.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff2705860e7492cff4c6fa7ace9133cb47bc0bb9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/878e5e59-2378-49c6-bcb0-0a60edbf870d.93e99953d7ab874464599f6eb9ece0ac00126114.zh-cn.xlf",
    $null,
    $null,
    "878e5e59-2378-49c6-bcb0-0a60edbf870d.93e99953d7ab874464599f6eb9ece0ac00126114.zh-cn.xlf")

$wsZh.Range("G2").Value = "2016-01-11 13:35:11"
$wsZh.Range("G3").Value = "2016-01-11 13:35:11"

$wsZh.Range("H2").Value = "Include"
$wsZh.Range("H3").Value = "Include"

# --- de-de sheet ------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

$wsDe.Hyperlinks.Add(
    $wsDe.Range("E2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d0f4f5c5fc27c12fe0d02e1e3d4264163c74298/e2e/6a104238-6086-4051-b5c2-9444b42ccc29.md",
    $null,
    $null,
    "6a104238-6086-4051-b5c2-9444b42ccc29.md")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5bad1f9fb6d19b225ebd1d2b7b1f9013e32c8dd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6a104238-6086-4051-b5c2-9444b42ccc29.4a63ded2fe0ab905558bc9b658fa27f0c8174700.de-de.xlf",
    $null,
    $null,
    "6a104238-6086-4051-b5c2-9444b42ccc29.4a63ded2fe0ab905558bc9b658fa27f0c8174700.de-de.xlf")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("E3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d0f4f5c5fc27c12fe0d02e1e3d4264163c74298/e2e/878e5e59-2378-49c6-bcb0-0a60edbf870d.md",
    $null,
    $null,
    "878e5e59-2378-49c6-bcb0-0a60edbf870d.md")

$wsDe.Hyperlinks.Add(
    $wsDe.Range("F3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f5bad1f9fb6d19b225ebd1d2b7b1f9013e32c8dd/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/878e5e59-2378-49c6-bcb0-0a60edbf870d.93e99953d7ab874464599f6eb9ece0ac00126114.de-de.xlf",
    $null,
    $null,
    "878e5e59-2378-49c6-bcb0-0a60edbf870d.93e99953d7ab874464599f6eb9ece0ac00126114.de-de.xlf")

$wsDe.Range("G2").Value = "2016-01-11 13:35:40"
$wsDe.Range("G3").Value = "2016-01-11 13:35:40"

$wsDe.Range("H2").Value = "Include"
$wsDe.Range("H3").Value = "Include"

Write-Output "Handback report generated."
